$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Switch example group names (formerly Alice/Bob/Claire/David/Elaine)
$ws.Range("B12").Value = "Veselin"
$ws.Range("B13").Value = "Rawda"
$ws.Range("B14").Value = "Hannah"
$ws.Range("B15").Value = "Mirit"
$ws.Range("B16").Value = "Bogdana"

# New row for an additional team member
$ws.Range("B17").Value = "Martin"
$ws.Range("E17").Value = 1

# Update the active selection
$ws.Activate()
$ws.Range("F13").Select()
